{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst anchorText = \"You can also switch/move to a particular desktop by number. The default is to use the plain number keys. You can switch to the F1-12 keys but these are very commonly already assigned for other uses.\";\n\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === anchorText) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Anchor paragraph not found\");\n}\n\n// Insert the body paragraph right after the anchor paragraph first (so it\n// inherits the anchor's implicit \"Normal\" style, i.e. no explicit pStyle),\n// then insert the heading paragraph immediately before it and mark it as\n// Heading2. This order keeps the new body paragraph free of an explicit\n// style reference, matching a plain Normal paragraph.\nconst bodyParagraph = anchor.insertParagraph(\n  \"Virtual Desktop Grid Switcher fixes an issue in Windows 10 where switching from a desktop which is empty to a desktop which had an activate window does not reactivate that window.\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\nconst headingParagraph = bodyParagraph.insertParagraph(\n  \"Window Activation on Switch From Empty Desktop\",\n  Word.InsertLocation.before\n);\nheadingParagraph.styleBuiltIn = Word.BuiltInStyleName.heading2;\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the paragraph that ends the \"Key Assignment\" section, right before\n# the \"Default Browser Activation\" heading - this is where the two new\n# paragraphs need to be inserted.\n$anchorText = \"You can also switch/move to a particular desktop by number. The default is to use the plain number keys. You can switch to the F1-12 keys but these are very commonly already assigned for other uses.\"\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd(\"`r\", \"`a\") -eq $anchorText) {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Anchor paragraph not found\"\n}\n\n# Insert a new (as yet empty) paragraph right after the anchor paragraph -\n# it inherits the anchor's implicit \"Normal\" style.\n$r = $target.Range\n$r.Collapse(0)   # wdCollapseEnd\n$r.InsertParagraphAfter()\n\n$bodyPara = $target.Next()\n$bodyPara.Range.Text = \"Virtual Desktop Grid Switcher fixes an issue in Windows 10 where switching from a desktop which is empty to a desktop which had an activate window does not reactivate that window.\"\n\n# Insert the new heading paragraph before the body paragraph and mark it\n# as a Heading 2, so the new body paragraph keeps its implicit Normal style.\n$br = $bodyPara.Range\n$br.Collapse(1)  # wdCollapseStart\n$br.InsertParagraphBefore()\n\n$headingPara = $target.Next()\n$headingPara.Range.Text = \"Window Activation on Switch From Empty Desktop\"\n$headingPara.Style = \"Heading 2\"\n"}
